$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '34.683.77'
$ws.Range("E2").Value = '  -1.99%  '
$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '1.807.52'
$ws.Range("E3").Value = '  -1.89%  '
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '232.33'
$ws.Range("E5").Value = '  +1.04%  '
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '0.602'
$ws.Range("E6").Value = '  -1.19%  '
$ws.Range("E7").Value = '  +0.34%  '
$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '39.33'
$ws.Range("E8").Value = '  -8.26%  '
$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.323'
$ws.Range("E9").Value = '  +5.03%  '
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.0680'
$ws.Range("E10").Value = '  -1.97%  '
$ws.Range("E11").Value = '  -1.85%  '
$ws.Range("E12").Value = '  -1.77%  '
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '1.820.80'
$ws.Range("E13").Value = '  -1.01%  '
$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '0.667'
$ws.Range("E14").Value = '  -0.51%  '
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '11.00'
$ws.Range("E15").Value = '  -2.48%  '
$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '4.56'
$ws.Range("E16").Value = '  -2.08%  '
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '34.686.66'
$ws.Range("E17").Value = '  -1.93%  '
$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '69.44'
$ws.Range("E18").Value = '  -1.04%  '
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '0.0₃0784'
$ws.Range("E19").Value = '  -1.27%  '
$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '239.29'
$ws.Range("E20").Value = '  -2.13%  '
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '11.91'
$ws.Range("E21").Value = '  -1.37%  '
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '4.62'
$ws.Range("E22").Value = '  -0.43%  '
$ws.Range("E23").Value = '  +0.22%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '2.24'
$ws.Range("E24").Value = '  +2.12%  '
$ws.Range("E25").Value = '  +1.80%  '
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '7.70'
$ws.Range("E26").Value = '  -2.51%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '17.15'
$ws.Range("E27").Value = '  -3.39%  '
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '0.119'
$ws.Range("E28").Value = '  -1.84%  '
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '1.54'
$ws.Range("E29").Value = '  +10.86%  '
$ws.Range("E30").Value = '  +0.33%  '
$ws.Range("E31").Value = '  +1.98%  '
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '0.0545'
$ws.Range("E32").Value = '  +0.07%  '
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '3.96'
$ws.Range("E33").Value = '  -2.66%  '
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '1.30'
$ws.Range("E34").Value = '  +19.52%  '
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '1.77'
$ws.Range("E35").Value = '  -4.44%  '
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '0.699'
$ws.Range("E36").Value = '  +1.94%  '
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '91.36'
$ws.Range("E37").Value = '  -4.42%  '
$ws.Range("E38").Value = '  +5.52%  '
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '1.316.61'
$ws.Range("E39").Value = '  -1.99%  '
$ws.Range("E40").Value = '  -1.16%  '
$ws.Range("E41").Value = '  +0.41%  '
$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '0.961'
$ws.Range("E42").Value = '  -3.84%  '
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '14.28'
$ws.Range("E43").Value = '  -3.42%  '
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '2.21'
$ws.Range("E44").Value = '  -9.10%  '
$ws.Range("E45").Value = '  -5.34%  '
$ws.Range("E46").Value = '  -0.49%  '
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '0.0512'
$ws.Range("E47").Value = '  -1.50%  '
$ws.Range("E48").Value = '  -0.52%  '
$ws.Range("E49").Value = '  +0.27%  '
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '0.0667'
$ws.Range("E50").Value = '  +7.27%  '
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '98.48'
$ws.Range("E51").Value = '  -4.48%  '
